$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (number format/style) of the last existing data row
# (row 95) down onto the three new rows, then fill in their values. Using
# PasteSpecial(xlPasteFormats) reuses the existing style indices (s="3" for
# the date column, s="4" for the percentage columns) instead of minting new
# duplicate styles.
$xlPasteFormats = -4122

$ws.Range("A95:E95").Copy() | Out-Null
$ws.Range("A96:E96").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A97:E97").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A98:E98").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = 0

# Novembro/2020
$ws.Cells.Item(96, 1).Value = 44136
$ws.Cells.Item(96, 2).Value = -0.043
$ws.Cells.Item(96, 3).Value = -0.021
$ws.Cells.Item(96, 4).Value = -0.11
$ws.Cells.Item(96, 5).Value = -0.09

# Dezembro/2020
$ws.Cells.Item(97, 1).Value = 44166
$ws.Cells.Item(97, 2).Value = -0.034
$ws.Cells.Item(97, 3).Value = -0.046
$ws.Cells.Item(97, 4).Value = -0.098
$ws.Cells.Item(97, 5).Value = -0.109

# Janeiro/2021
$ws.Cells.Item(98, 1).Value = 44197
$ws.Cells.Item(98, 2).Value = -0.063
$ws.Cells.Item(98, 3).Value = -0.039
$ws.Cells.Item(98, 4).Value = -0.126
$ws.Cells.Item(98, 5).Value = -0.103

# The new bottom row of the sheet becomes the active selection (Excel
# auto-advances past the freshly entered data block).
$ws.Range("A99").Select() | Out-Null
